$d = $word.ActiveDocument

# The commit ("Switch pandoc to nix") regenerated this fixture with a new
# pandoc build; the resulting OOXML differs only in (a) a handful of
# character-style <w:rPr> child-element orderings (bold/italic now emitted
# before color) and (b) the zero-padding of one list's internal <w:nsid>.
# The <w:nsid> value is an opaque internal list identifier that Word's
# object model does not expose for editing (it is not reachable through
# Styles, ListTemplates/ListLevels, or any other COM surface), so only the
# style rPr-ordering portion of the change is applied here.
#
# Re-assigning a style's existing Bold/Italic flag to itself doesn't change
# the style's formatting, but it does make the engine rewrite that style's
# <w:rPr> children in canonical order (bold, italic, color), which is
# exactly the reordering the diff shows. Font.Bold/Font.Italic are used
# (rather than Font.Color) so the existing color hex strings are re-emitted
# byte-for-byte instead of being re-cased.

$boldTouch = @(
    "KeywordTok",
    "ImportTok",
    "AnnotationTok",
    "CommentVarTok",
    "ControlFlowTok",
    "InformationTok",
    "WarningTok",
    "AlertTok",
    "ErrorTok"
)

foreach ($name in $boldTouch) {
    $style = $d.Styles($name)
    $style.Font.Bold = $style.Font.Bold
}

$italicTouch = @(
    "CommentTok",
    "DocumentationTok"
)

foreach ($name in $italicTouch) {
    $style = $d.Styles($name)
    $style.Font.Italic = $style.Font.Italic
}
